$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 4715.0586
$ws.Range("I6").Value = 6373.4165
$ws.Range("K6").Value = 19120.2495
$ws.Range("M6").Value = -19008.2495
$ws.Range("H112").Value = 5161.643
$ws.Range("J112").Value = 4264.28
$ws.Range("L112").Value = 12792.84
$ws.Range("N112").Value = -15008.84
$ws.Range("H113").Value = 3200
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3200
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9708
$ws.Range("H137").Value = 15387.34
$ws.Range("I137").Value = 8730.691999999999
$ws.Range("J137").Value = 17726.162
$ws.Range("K137").Value = 26192.076
$ws.Range("L137").Value = 53178.486
$ws.Range("M137").Value = -23642.076
$ws.Range("N137").Value = -58278.486

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 223038.44
$ws.Range("I13").Value = 223038.44
$ws.Range("K13").Value = 223038.44
$ws.Range("M13").Value = -222894.44
$ws.Range("H32").Value = 2204600.5
$ws.Range("I32").Value = 3237794.8
$ws.Range("K32").Value = 3237794.8
$ws.Range("M32").Value = -3237507.8
$ws.Range("H45").Value = 3500.625
$ws.Range("I45").Value = 1667.5
$ws.Range("J45").Value = 9000
$ws.Range("K45").Value = 1667.5
$ws.Range("L45").Value = 9000
$ws.Range("M45").Value = -1290.5
$ws.Range("N45").Value = -9754
$ws.Range("H74").Value = 19751.75
$ws.Range("I74").Value = 2705.9092
$ws.Range("J74").Value = 34175.152
$ws.Range("K74").Value = 2705.9092
$ws.Range("L74").Value = 34175.152
$ws.Range("M74").Value = -1831.9092
$ws.Range("N74").Value = -35923.152
$ws.Range("H77").Value = 19751.75
$ws.Range("I77").Value = 2705.9092
$ws.Range("J77").Value = 34175.152
$ws.Range("K77").Value = 13529.546
$ws.Range("L77").Value = 170875.76
$ws.Range("M77").Value = -9161.546
$ws.Range("N77").Value = -179611.76

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 24458.242
$ws.Range("I20").Value = 6627
$ws.Range("J20").Value = 41240.59
$ws.Range("K20").Value = 6627
$ws.Range("L20").Value = 41240.59
$ws.Range("M20").Value = -6380
$ws.Range("N20").Value = -41734.59
$ws.Range("H24").Value = 5666.6665
$ws.Range("J24").Value = 6500
$ws.Range("L24").Value = 6500
$ws.Range("N24").Value = -6970
$ws.Range("H86").Value = 4767.357
$ws.Range("I86").Value = 4941.769
$ws.Range("K86").Value = 4941.769
$ws.Range("M86").Value = -3818.769
$ws.Range("H89").Value = 4767.357
$ws.Range("I89").Value = 4941.769
$ws.Range("K89").Value = 24708.845
$ws.Range("M89").Value = -19092.845
$ws.Range("H134").Value = 11602.744
$ws.Range("I134").Value = 7291.75
$ws.Range("K134").Value = 21875.25
$ws.Range("M134").Value = -19340.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30124.238
$ws.Range("I31").Value = 19182.695
$ws.Range("J31").Value = 43369.26
$ws.Range("K31").Value = 19182.695
$ws.Range("L31").Value = 43369.26
$ws.Range("M31").Value = -18887.695
$ws.Range("N31").Value = -43959.26
$ws.Range("H34").Value = 30124.238
$ws.Range("I34").Value = 19182.695
$ws.Range("J34").Value = 43369.26
$ws.Range("K34").Value = 19182.695
$ws.Range("L34").Value = 43369.26
$ws.Range("M34").Value = -18980.695
$ws.Range("N34").Value = -43773.26
$ws.Range("H62").Value = 5697.1665
$ws.Range("I62").Value = 5086.6
$ws.Range("J62").Value = 8750
$ws.Range("K62").Value = 5086.6
$ws.Range("L62").Value = 8750
$ws.Range("M62").Value = -4462.6
$ws.Range("N62").Value = -9998
$ws.Range("H65").Value = 5697.1665
$ws.Range("I65").Value = 5086.6
$ws.Range("J65").Value = 8750
$ws.Range("K65").Value = 25433
$ws.Range("L65").Value = 43750
$ws.Range("M65").Value = -22313
$ws.Range("N65").Value = -49990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 21318056
$ws.Range("I4").Value = 31626462
$ws.Range("K4").Value = 94879386
$ws.Range("M4").Value = -94879274
$ws.Range("H40").Value = 432.6316
$ws.Range("I40").Value = 435.29413
$ws.Range("K40").Value = 1741.17652
$ws.Range("M40").Value = -1672.17652

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 17812.285
$ws.Range("I126").Value = 20349.5
$ws.Range("J126").Value = 16797.4
$ws.Range("K126").Value = 61048.5
$ws.Range("L126").Value = 50392.2
$ws.Range("M126").Value = -58578.5
$ws.Range("N126").Value = -55332.2
$ws.Range("H132").Value = 9201.143
$ws.Range("I132").Value = 10774.272
$ws.Range("K132").Value = 32322.816
$ws.Range("M132").Value = -29792.816

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11549.571
$ws.Range("I7").Value = 7266.6665
$ws.Range("K7").Value = 7266.6665
$ws.Range("M7").Value = -7154.6665
$ws.Range("H40").Value = 7436.6206
$ws.Range("I40").Value = 4137.353
$ws.Range("K40").Value = 4137.353
$ws.Range("M40").Value = -4001.353
$ws.Range("H46").Value = 2544.5334
$ws.Range("I46").Value = 2095.7144
$ws.Range("J46").Value = 2937.25
$ws.Range("K46").Value = 2095.7144
$ws.Range("L46").Value = 2937.25
$ws.Range("M46").Value = -1907.7144
$ws.Range("N46").Value = -3313.25
$ws.Range("H68").Value = 10866.167
$ws.Range("J68").Value = 26999.5
$ws.Range("L68").Value = 26999.5
$ws.Range("N68").Value = -28497.5
$ws.Range("H71").Value = 10866.167
$ws.Range("J71").Value = 26999.5
$ws.Range("L71").Value = 134997.5
$ws.Range("N71").Value = -142485.5
$ws.Range("H126").Value = 11549.571
$ws.Range("I126").Value = 7266.6665
$ws.Range("K126").Value = 21799.9995
$ws.Range("M126").Value = -19329.9995
$ws.Range("H136").Value = 22410.822
$ws.Range("I136").Value = 19093.438
$ws.Range("J136").Value = 26834
$ws.Range("K136").Value = 57280.314
$ws.Range("L136").Value = 80502
$ws.Range("M136").Value = -54730.314
$ws.Range("N136").Value = -85602

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 64833
$ws.Range("I126").Value = 96712.25
$ws.Range("J126").Value = 1074.5
$ws.Range("K126").Value = 290136.75
$ws.Range("L126").Value = 3223.5
$ws.Range("M126").Value = -287666.75
$ws.Range("N126").Value = -8163.5
